$d = $word.ActiveDocument

# --- Update the title paragraph ---
$d.Paragraphs.Item(1).Range.Text = "數學 - 應用題 - 典型應用題 - 盈虧問題 - 住宿、乘車、乘船"

# --- Update existing question paragraphs (2) through (7) i.e. questions 1-6 ---
$d.Paragraphs.Item(2).Range.Text = "(1) 希望小學學生乘遊覽車到西子灣郊遊。如果每車坐65人，則有15人不能乘車；如果每車多坐5人，恰好多餘了一輛車。問：共有【17】輛汽車、有【1120】學生"
$d.Paragraphs.Item(3).Range.Text = "(2) 某校安排學生宿舍，如果每間住5人則有14人沒有床位；如果每間住7人，則多出4個床位。問：宿舍【9】間、住宿生【59】人。"
$d.Paragraphs.Item(4).Range.Text = "(3) 同學們去划船，如果每條船坐5人，則少兩條船，如果每條船坐7人，則多出兩條船，問：共有【12】條船、有【70】個同學。"
$d.Paragraphs.Item(5).Range.Text = "(4) 某學校有一些學生住校，每間宿舍住10人，空出床位24張，如果每間宿舍住8人，則空出床位2張。問：學校共有【11】間宿舍、住宿學生有【86】人。"
$d.Paragraphs.Item(6).Range.Text = "(5) 工廠新建一宿舍，每間住4人，則有34人沒床位，每間住6人，則又多5間房，問：共有【162】名工人要安排住宿。"
$d.Paragraphs.Item(7).Range.Text = "(6) 消防局調動一批消防員分乘一批車輛趕往災區救災。原計劃每輛汽車乘32人，則多出5人，他們被安排乘坐在其中的某輛車上，行進中由於緊急狀況調走一輛車，這時只好重新安排每輛車乘35人，這樣多出7人，他們被安排在其中的某輛車上，問：原來共有【11】輛車、共派出【357】名消防員"

# --- Append new question paragraphs (7) through (10) ---
$newTexts = @(
  "(7) 學校給一些新同學分配宿舍，如果每個房間住12人，則34人沒有位置，如果每個房間住14人，則空出4個房間。問：宿舍有【45】間、有【574】學生",
  "(8) 某校參加六一杯小學數學競賽，原定考場若干個。如果增加2個考場，每個考場正好坐24人；如果減少2個考場，每個考場正好坐30人。問：參加這次競賽的學生共有【480】人、原定考場【18】個 ",
  "(9) 學校給新生分配宿舍，如果每間住8人，則少2間房；如果每間住10人，則多出2間房。問：共有【18】間房分給新生、新生有【160】人",
  "(10) 同學們去划船，如果每艘船坐4人，則少3艘船；若每艘船坐6人，還有2人留在岸邊。問：有【32】名同學去划船，共租【5】艘船。"
)

foreach ($t in $newTexts) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = $t
}
